$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "29.381.63"
$ws.Range("E2").Value = "  +0.07%  "

Set-TextValue "D3" "1.841.50"
$ws.Range("E3").Value = "  -0.19%  "

Set-TextValue "D4" "0.9992"
$ws.Range("E4").Value = "  +0.18%  "

Set-TextValue "D5" "239.15"
$ws.Range("E5").Value = "  -0.38%  "

Set-TextValue "D6" "0.6266"
$ws.Range("E6").Value = "  -0.14%  "

Set-TextValue "D7" "1.001"

Set-TextValue "D8" "0.07419"
$ws.Range("E8").Value = "  -0.98%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue "D9" "0.2891"
$ws.Range("E9").Value = "  -0.29%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D10" "24.93"
$ws.Range("E10").Value = "  +1.87%  "

Set-TextValue "D11" "0.07718"
$ws.Range("E11").Value = "  -0.21%  "

Set-TextValue "D12" "1.831.40"
$ws.Range("E12").Value = "  -0.73%  "

Set-TextValue "D13" "4.977"
$ws.Range("E13").Value = "  -0.38%  "

Set-TextValue "D14" "0.6739"
$ws.Range("E14").Value = "  -1.03%  "

Set-TextValue "D15" "0.00001030"
$ws.Range("E15").Value = "  -1.80%  "

Set-TextValue "D16" "81.74"
$ws.Range("E16").Value = "  -0.41%  "

Set-TextValue "D17" "6.207"
$ws.Range("E17").Value = "  +0.25%  "

Set-TextValue "D18" "29.410.16"
$ws.Range("E18").Value = "  +0.09%  "

Set-TextValue "D19" "233.76"
$ws.Range("E19").Value = "  +2.08%  "

Set-TextValue "D20" "12.31"
$ws.Range("E20").Value = "  -0.19%  "

Set-TextValue "D21" "1.001"
$ws.Range("E21").Value = "  +0.20%  "

Set-TextValue "D22" "7.290"
$ws.Range("E22").Value = "  -2.74%  "

Set-TextValue "D23" "1.001"
$ws.Range("E23").Value = "  +0.25%  "

Set-TextValue "D24" "158.64"
$ws.Range("E24").Value = "  +0.04%  "

Set-TextValue "D25" "8.506"
$ws.Range("E25").Value = "  +0.95%  "

$ws.Range("E26").Value = "  -1.90%  "

Set-TextValue "D27" "17.29"
$ws.Range("E27").Value = "  -1.22%  "

Set-TextValue "D28" "0.07310"
$ws.Range("E28").Value = "  +12.61%  "

Set-TextValue "D29" "1.466"
$ws.Range("E29").Value = "  +4.01%  "

Set-TextValue "D30" "1.479"
$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D31" "4.041"
$ws.Range("E31").Value = "  -1.17%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D32" "4.029"
$ws.Range("E32").Value = "  -1.58%  "

Set-TextValue "D33" "1.817"
$ws.Range("E33").Value = "  -0.86%  "

Set-TextValue "D34" "1.139"
$ws.Range("E34").Value = "  -0.26%  "

Set-TextValue "D35" "0.6964"
$ws.Range("E35").Value = "  -0.19%  "

Set-TextValue "D36" "2.576"
$ws.Range("E36").Value = "  -0.12%  "

Set-TextValue "D37" "0.01842"
$ws.Range("E37").Value = "  +0.26%  "

Set-TextValue "D38" "6.912"
$ws.Range("E38").Value = "  +2.24%  "

$ws.Range("E39").Value = "  -0.85%  "

Set-TextValue "D40" "1.233.15"
$ws.Range("E40").Value = "  -2.79%  "

Set-TextValue "D41" "0.9556"
$ws.Range("E41").Value = "  +4.14%  "

$ws.Range("E42").Value = "  +0.22%  "

Set-TextValue "D43" "1.997.04"
$ws.Range("E43").Value = "  -0.57%  "

Set-TextValue "D44" "100.89"
$ws.Range("E44").Value = "  -0.41%  "

Set-TextValue "D45" "65.36"
$ws.Range("E45").Value = "  -1.30%  "

$ws.Range("E46").Value = "  +1.88%  "

Set-TextValue "D47" "1.715"
$ws.Range("E47").Value = "  -0.68%  "

Set-TextValue "D48" "6.947"
$ws.Range("E48").Value = "  -1.92%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "8.854"
$ws.Range("E49").Value = "  -1.78%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.1133"
$ws.Range("E50").Value = "  -2.87%  "

Set-TextValue "D51" "0.3896"
$ws.Range("E51").Value = "  -1.76%  "
